# empleados_bac.xlsx update: switch the "nit" column header to "id" and add a
# new employee record (BRENDA GRISELDA ROMERO HERNANDEZ) as row 3, pushing the
# remaining employees down by one row. The new row's id cell is highlighted
# (pale-yellow fill, medium grey border, Calibri font, left aligned + wrapped,
# integer number format) to flag that it uses a short DUI-style id rather than
# the usual 13-digit NIT used by the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cell C1: "nit" -> "id"
$ws.Range("C1").Value = "id"

# 2. Insert a new blank row at position 3; existing rows 3..29 shift to 4..30.
#    (Insert() copies the formatting of the row above by default.)
$ws.Rows.Item(3).Insert() | Out-Null

# 3. The new row's alias/nombre/cuenta cells should look like plain,
#    unformatted cells (same as the header row) rather than inheriting the
#    data-row styling that Insert() copied down - reset them via a
#    formats-only paste from the header row.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A3:B3").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 4. Fill in the new employee's data.
$ws.Range("A3").Value = "BRENDA GRISELDA ROMERO HERNANDEZ"
$ws.Range("B3").Value = "BRENDA GRISELDA ROMERO HERNANDEZ"
$ws.Range("C3").Value = 55544433
$ws.Range("D3").Value = 122641731

# 5. Highlight the new "id" cell (C3) to call out that it's a DUI, not a NIT.
$c3 = $ws.Range("C3")
$c3.NumberFormat = "0"
$c3.Font.Name = "Calibri"
$c3.Font.Size = 11
$c3.Font.Color = 0
$c3.Interior.Pattern = 1
$c3.Interior.Color = 16777215
$c3.Interior.PatternColor = 13434879
$c3.Borders.Color = 13948116
$c3.Borders.Weight = -4138
$c3.HorizontalAlignment = -4131
$c3.WrapText = $true

# 6. The new row is a touch taller than the rest.
$ws.Rows.Item(3).RowHeight = 13.8

# 7. Leave the selection on the newly highlighted id cell.
$ws.Range("C4").Select() | Out-Null
